# The "Gaz" sheet gains a new column E with a gas price (€/MWh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# Copy the header formatting from D1 (bold, bordered, centered) onto the
# new header cell E1, then set its text (xlPasteFormats = -4122).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E1").Value = "Prix Gaz (€/MWh)"

# Fill in the gas price values for each of the three data rows.
$ws.Range("E2").Value = 37.15
$ws.Range("E3").Value = 38.95
$ws.Range("E4").Value = 38.45
